$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "Additional ABAQUS Resources"
$ws.Range("B2").Value = "Resources"
$wb.Worksheets.Item("Material Properties").Activate()
